$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# AF5: related_works value changed for row 5
$ws.Range("AF5").Value = 'c("https://openalex.org/W2090007074", "https://openalex.org/W4281261429", "https://openalex.org/W4205683907", "https://openalex.org/W4327921319", "https://openalex.org/W1976433721", "https://openalex.org/W2414935542", "https://openalex.org/W4387268337", "https://openalex.org/W2438041563", "https://openalex.org/W2347635326", "https://openalex.org/W2561050005")'

# F6, G6, H6: source name, source id, host organization for row 6 (was N/A, now filled in)
$ws.Range("F6").Value = "medRxiv (Cold Spring Harbor Laboratory)"
$ws.Range("G6").Value = "https://openalex.org/S4306400573"
$ws.Range("H6").Value = "Cold Spring Harbor Laboratory"
